$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap F:V contents between paired rows (matches swapped betexplorer
# fixture order for rows 5/6, 17/18, 35/36, 41/42, 49/50) ---

# Rows 5 <-> 6
$ws.Cells.Item(5, 6).Value = "Magra"
$ws.Cells.Item(5, 7).Value = 0
$ws.Cells.Item(5, 8).Value = "Kabylie"
$ws.Cells.Item(5, 9).Value = 1
$ws.Cells.Item(5, 10).Value = 2.84
$ws.Cells.Item(5, 11).Value = "15/09/2023 13:42"
$ws.Cells.Item(5, 12).Value = 3.8
$ws.Cells.Item(5, 13).Value = "16/09/2023 16:12"
$ws.Cells.Item(5, 14).Value = 2.63
$ws.Cells.Item(5, 15).Value = "15/09/2023 13:42"
$ws.Cells.Item(5, 16).Value = 2.84
$ws.Cells.Item(5, 17).Value = "16/09/2023 15:03"
$ws.Cells.Item(5, 18).Value = 2.72
$ws.Cells.Item(5, 19).Value = "15/09/2023 13:42"
$ws.Cells.Item(5, 20).Value = 2.26
$ws.Cells.Item(5, 21).Value = "16/09/2023 16:12"
$ws.Cells.Item(5, 22).Value = "https://www.betexplorer.com/football/algeria/ligue-1/magra-kabylie/YFXa8c8H/"
$ws.Cells.Item(6, 6).Value = "MC Alger"
$ws.Cells.Item(6, 7).Value = 4
$ws.Cells.Item(6, 8).Value = "Ben Aknoun"
$ws.Cells.Item(6, 9).Value = 0
$ws.Cells.Item(6, 10).Value = 1.4
$ws.Cells.Item(6, 11).Value = "16/09/2023 03:43"
$ws.Cells.Item(6, 12).Value = 1.3
$ws.Cells.Item(6, 13).Value = "16/09/2023 10:40"
$ws.Cells.Item(6, 14).Value = 4.19
$ws.Cells.Item(6, 15).Value = "16/09/2023 03:43"
$ws.Cells.Item(6, 16).Value = 4.81
$ws.Cells.Item(6, 17).Value = "16/09/2023 16:47"
$ws.Cells.Item(6, 18).Value = 8.529999999999999
$ws.Cells.Item(6, 19).Value = "16/09/2023 03:43"
$ws.Cells.Item(6, 20).Value = 12.64
$ws.Cells.Item(6, 21).Value = "16/09/2023 16:47"
$ws.Cells.Item(6, 22).Value = "https://www.betexplorer.com/football/algeria/ligue-1/mc-alger-es-ben-aknoun/WjyqCu9h/"

# Rows 17 <-> 18
$ws.Cells.Item(17, 6).Value = "US Souf"
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = "Oran"
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 2.49
$ws.Cells.Item(17, 11).Value = "28/09/2023 19:27"
$ws.Cells.Item(17, 12).Value = 2.14
$ws.Cells.Item(17, 13).Value = "29/09/2023 13:29"
$ws.Cells.Item(17, 14).Value = 2.88
$ws.Cells.Item(17, 15).Value = "28/09/2023 19:27"
$ws.Cells.Item(17, 16).Value = 2.74
$ws.Cells.Item(17, 17).Value = "29/09/2023 14:49"
$ws.Cells.Item(17, 18).Value = 3.18
$ws.Cells.Item(17, 19).Value = "28/09/2023 19:27"
$ws.Cells.Item(17, 20).Value = 4.43
$ws.Cells.Item(17, 21).Value = "29/09/2023 15:47"
$ws.Cells.Item(17, 22).Value = "https://www.betexplorer.com/football/algeria/ligue-1/us-souf-oran/6qOsFaSf/"
$ws.Cells.Item(18, 6).Value = "Khenchela"
$ws.Cells.Item(18, 7).Value = 2
$ws.Cells.Item(18, 8).Value = "Kabylie"
$ws.Cells.Item(18, 9).Value = 1
$ws.Cells.Item(18, 10).Value = 2.63
$ws.Cells.Item(18, 11).Value = "28/09/2023 04:12"
$ws.Cells.Item(18, 12).Value = 2.05
$ws.Cells.Item(18, 13).Value = "29/09/2023 16:41"
$ws.Cells.Item(18, 14).Value = 2.62
$ws.Cells.Item(18, 15).Value = "28/09/2023 04:12"
$ws.Cells.Item(18, 16).Value = 2.75
$ws.Cells.Item(18, 17).Value = "29/09/2023 16:41"
$ws.Cells.Item(18, 18).Value = 3.02
$ws.Cells.Item(18, 19).Value = "28/09/2023 04:12"
$ws.Cells.Item(18, 20).Value = 4.88
$ws.Cells.Item(18, 21).Value = "29/09/2023 16:27"
$ws.Cells.Item(18, 22).Value = "https://www.betexplorer.com/football/algeria/ligue-1/khenchela-kabylie/pUZYGLcr/"

# Rows 35 <-> 36
$ws.Cells.Item(35, 6).Value = "Magra"
$ws.Cells.Item(35, 7).Value = 3
$ws.Cells.Item(35, 8).Value = "Ben Aknoun"
$ws.Cells.Item(35, 9).Value = 1
$ws.Cells.Item(35, 10).Value = 1.61
$ws.Cells.Item(35, 11).Value = "11/11/2023 10:12"
$ws.Cells.Item(35, 12).Value = 1.62
$ws.Cells.Item(35, 13).Value = "11/11/2023 14:48"
$ws.Cells.Item(35, 14).Value = 3.51
$ws.Cells.Item(35, 15).Value = "11/11/2023 10:12"
$ws.Cells.Item(35, 16).Value = 3.56
$ws.Cells.Item(35, 17).Value = "11/11/2023 14:48"
$ws.Cells.Item(35, 18).Value = 5.97
$ws.Cells.Item(35, 19).Value = "11/11/2023 10:12"
$ws.Cells.Item(35, 20).Value = 6.3
$ws.Cells.Item(35, 21).Value = "11/11/2023 14:48"
$ws.Cells.Item(35, 22).Value = "https://www.betexplorer.com/football/algeria/ligue-1/magra-es-ben-aknoun/lCJE0FP6/"
$ws.Cells.Item(36, 6).Value = "Khenchela"
$ws.Cells.Item(36, 7).Value = 0
$ws.Cells.Item(36, 8).Value = "Biskra"
$ws.Cells.Item(36, 9).Value = 1
$ws.Cells.Item(36, 10).Value = 1.74
$ws.Cells.Item(36, 11).Value = "10/11/2023 03:13"
$ws.Cells.Item(36, 12).Value = 1.29
$ws.Cells.Item(36, 13).Value = "11/11/2023 10:24"
$ws.Cells.Item(36, 14).Value = 3.19
$ws.Cells.Item(36, 15).Value = "10/11/2023 03:13"
$ws.Cells.Item(36, 16).Value = 4.87
$ws.Cells.Item(36, 17).Value = "11/11/2023 14:54"
$ws.Cells.Item(36, 18).Value = 4.9
$ws.Cells.Item(36, 19).Value = "10/11/2023 03:13"
$ws.Cells.Item(36, 20).Value = 13.77
$ws.Cells.Item(36, 21).Value = "11/11/2023 14:54"
$ws.Cells.Item(36, 22).Value = "https://www.betexplorer.com/football/algeria/ligue-1/khenchela-biskra/GbL62yef/"

# Rows 41 <-> 42
$ws.Cells.Item(41, 6).Value = "Constantine"
$ws.Cells.Item(41, 7).Value = 0
$ws.Cells.Item(41, 8).Value = "Magra"
$ws.Cells.Item(41, 9).Value = 1
$ws.Cells.Item(41, 10).Value = 1.59
$ws.Cells.Item(41, 11).Value = "26/10/2023 04:42"
$ws.Cells.Item(41, 12).Value = 1.31
$ws.Cells.Item(41, 13).Value = "17/11/2023 15:24"
$ws.Cells.Item(41, 14).Value = 3.47
$ws.Cells.Item(41, 15).Value = "26/10/2023 04:42"
$ws.Cells.Item(41, 16).Value = 4.98
$ws.Cells.Item(41, 17).Value = "17/11/2023 15:24"
$ws.Cells.Item(41, 18).Value = 5.42
$ws.Cells.Item(41, 19).Value = "26/10/2023 04:42"
$ws.Cells.Item(41, 20).Value = 10.79
$ws.Cells.Item(41, 21).Value = "17/11/2023 15:24"
$ws.Cells.Item(41, 22).Value = "https://www.betexplorer.com/football/algeria/ligue-1/constantine-magra/jgeqGfPg/"
$ws.Cells.Item(42, 6).Value = "El Bayadh"
$ws.Cells.Item(42, 7).Value = 4
$ws.Cells.Item(42, 8).Value = "US Souf"
$ws.Cells.Item(42, 9).Value = 0
$ws.Cells.Item(42, 10).Value = 1.45
$ws.Cells.Item(42, 11).Value = "16/11/2023 03:42"
$ws.Cells.Item(42, 12).Value = 1.45
$ws.Cells.Item(42, 13).Value = "17/11/2023 15:18"
$ws.Cells.Item(42, 14).Value = 3.74
$ws.Cells.Item(42, 15).Value = "16/11/2023 03:42"
$ws.Cells.Item(42, 16).Value = 4.01
$ws.Cells.Item(42, 17).Value = "17/11/2023 15:18"
$ws.Cells.Item(42, 18).Value = 7.38
$ws.Cells.Item(42, 19).Value = "16/11/2023 03:42"
$ws.Cells.Item(42, 20).Value = 8.609999999999999
$ws.Cells.Item(42, 21).Value = "17/11/2023 15:18"
$ws.Cells.Item(42, 22).Value = "https://www.betexplorer.com/football/algeria/ligue-1/el-bayadh-us-souf/UqfmFEv0/"

# Rows 49 <-> 50
$ws.Cells.Item(49, 6).Value = "Magra"
$ws.Cells.Item(49, 7).Value = 1
$ws.Cells.Item(49, 8).Value = "El Bayadh"
$ws.Cells.Item(49, 9).Value = 1
$ws.Cells.Item(49, 10).Value = 2.07
$ws.Cells.Item(49, 11).Value = "23/11/2023 07:47"
$ws.Cells.Item(49, 12).Value = 2.15
$ws.Cells.Item(49, 13).Value = "24/11/2023 14:19"
$ws.Cells.Item(49, 14).Value = 2.88
$ws.Cells.Item(49, 15).Value = "23/11/2023 07:47"
$ws.Cells.Item(49, 16).Value = 2.9
$ws.Cells.Item(49, 17).Value = "24/11/2023 15:12"
$ws.Cells.Item(49, 18).Value = 4.25
$ws.Cells.Item(49, 19).Value = "23/11/2023 07:47"
$ws.Cells.Item(49, 20).Value = 4.06
$ws.Cells.Item(49, 21).Value = "24/11/2023 14:19"
$ws.Cells.Item(49, 22).Value = "https://www.betexplorer.com/football/algeria/ligue-1/magra-el-bayadh/jog1n073/"
$ws.Cells.Item(50, 6).Value = "Paradou"
$ws.Cells.Item(50, 7).Value = 0
$ws.Cells.Item(50, 8).Value = "Constantine"
$ws.Cells.Item(50, 9).Value = 0
$ws.Cells.Item(50, 10).Value = 1.79
$ws.Cells.Item(50, 11).Value = "23/11/2023 07:47"
$ws.Cells.Item(50, 12).Value = 1.98
$ws.Cells.Item(50, 13).Value = "24/11/2023 15:00"
$ws.Cells.Item(50, 14).Value = 3.27
$ws.Cells.Item(50, 15).Value = "23/11/2023 07:47"
$ws.Cells.Item(50, 16).Value = 3.08
$ws.Cells.Item(50, 17).Value = "24/11/2023 15:00"
$ws.Cells.Item(50, 18).Value = 4.96
$ws.Cells.Item(50, 19).Value = "23/11/2023 07:47"
$ws.Cells.Item(50, 20).Value = 4.37
$ws.Cells.Item(50, 21).Value = "24/11/2023 15:00"
$ws.Cells.Item(50, 22).Value = "https://www.betexplorer.com/football/algeria/ligue-1/paradou-constantine/dOgcmthc/"

# --- Append new match rows 67-70 ---
$ws.Range("A2").Copy()
$ws.Range("A67:A70").PasteSpecial(-4122)
$ws.Range("E2").Copy()
$ws.Range("E67:E70").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 67
$ws.Cells.Item(67, 1).Value = 66
$ws.Cells.Item(67, 2).Value = "algeria"
$ws.Cells.Item(67, 3).Value = "ligue-1"
$ws.Cells.Item(67, 4).Value = "2023-2024"
$ws.Cells.Item(67, 5).Value = 45275.63541666666
$ws.Cells.Item(67, 6).Value = "El Bayadh"
$ws.Cells.Item(67, 7).Value = 1
$ws.Cells.Item(67, 8).Value = "USM Alger"
$ws.Cells.Item(67, 9).Value = 1
$ws.Cells.Item(67, 10).Value = 2.12
$ws.Cells.Item(67, 11).Value = "14/12/2023 01:12"
$ws.Cells.Item(67, 12).Value = 2.16
$ws.Cells.Item(67, 13).Value = "15/12/2023 15:13"
$ws.Cells.Item(67, 14).Value = 2.86
$ws.Cells.Item(67, 15).Value = "14/12/2023 01:12"
$ws.Cells.Item(67, 16).Value = 2.85
$ws.Cells.Item(67, 17).Value = "15/12/2023 15:10"
$ws.Cells.Item(67, 18).Value = 3.99
$ws.Cells.Item(67, 19).Value = "14/12/2023 01:12"
$ws.Cells.Item(67, 20).Value = 4.3
$ws.Cells.Item(67, 21).Value = "15/12/2023 15:13"
$ws.Cells.Item(67, 22).Value = "https://www.betexplorer.com/football/algeria/ligue-1/el-bayadh-usm-alger/Ywta5DYH/"

# Row 68
$ws.Cells.Item(68, 1).Value = 67
$ws.Cells.Item(68, 2).Value = "algeria"
$ws.Cells.Item(68, 3).Value = "ligue-1"
$ws.Cells.Item(68, 4).Value = "2023-2024"
$ws.Cells.Item(68, 5).Value = 45275.63541666666
$ws.Cells.Item(68, 6).Value = "Ben Aknoun"
$ws.Cells.Item(68, 7).Value = 1
$ws.Cells.Item(68, 8).Value = "CR Belouizdad"
$ws.Cells.Item(68, 9).Value = 1
$ws.Cells.Item(68, 10).Value = 6.41
$ws.Cells.Item(68, 11).Value = "14/12/2023 00:12"
$ws.Cells.Item(68, 12).Value = 5.81
$ws.Cells.Item(68, 13).Value = "15/12/2023 15:12"
$ws.Cells.Item(68, 14).Value = 3.37
$ws.Cells.Item(68, 15).Value = "14/12/2023 00:12"
$ws.Cells.Item(68, 16).Value = 3.59
$ws.Cells.Item(68, 17).Value = "15/12/2023 15:12"
$ws.Cells.Item(68, 18).Value = 1.56
$ws.Cells.Item(68, 19).Value = "14/12/2023 00:12"
$ws.Cells.Item(68, 20).Value = 1.65
$ws.Cells.Item(68, 21).Value = "15/12/2023 15:12"
$ws.Cells.Item(68, 22).Value = "https://www.betexplorer.com/football/algeria/ligue-1/es-ben-aknoun-cr-belouizdad/CjWzBH4t/"

# Row 69
$ws.Cells.Item(69, 1).Value = 68
$ws.Cells.Item(69, 2).Value = "algeria"
$ws.Cells.Item(69, 3).Value = "ligue-1"
$ws.Cells.Item(69, 4).Value = "2023-2024"
$ws.Cells.Item(69, 5).Value = 45275.69791666666
$ws.Cells.Item(69, 6).Value = "Constantine"
$ws.Cells.Item(69, 7).Value = 2
$ws.Cells.Item(69, 8).Value = "Kabylie"
$ws.Cells.Item(69, 9).Value = 0
$ws.Cells.Item(69, 10).Value = 1.79
$ws.Cells.Item(69, 11).Value = "14/12/2023 01:12"
$ws.Cells.Item(69, 12).Value = 1.71
$ws.Cells.Item(69, 13).Value = "15/12/2023 16:40"
$ws.Cells.Item(69, 14).Value = 3.01
$ws.Cells.Item(69, 15).Value = "14/12/2023 01:12"
$ws.Cells.Item(69, 16).Value = 3.25
$ws.Cells.Item(69, 17).Value = "15/12/2023 16:40"
$ws.Cells.Item(69, 18).Value = 4.95
$ws.Cells.Item(69, 19).Value = "14/12/2023 01:12"
$ws.Cells.Item(69, 20).Value = 6.09
$ws.Cells.Item(69, 21).Value = "15/12/2023 16:35"
$ws.Cells.Item(69, 22).Value = "https://www.betexplorer.com/football/algeria/ligue-1/constantine-kabylie/zVv9FcdP/"

# Row 70
$ws.Cells.Item(70, 1).Value = 69
$ws.Cells.Item(70, 2).Value = "algeria"
$ws.Cells.Item(70, 3).Value = "ligue-1"
$ws.Cells.Item(70, 4).Value = "2023-2024"
$ws.Cells.Item(70, 5).Value = 45275.75
$ws.Cells.Item(70, 6).Value = "ASO Chlef"
$ws.Cells.Item(70, 7).Value = 1
$ws.Cells.Item(70, 8).Value = "Paradou"
$ws.Cells.Item(70, 9).Value = 1
$ws.Cells.Item(70, 10).Value = 2
$ws.Cells.Item(70, 11).Value = "14/12/2023 01:12"
$ws.Cells.Item(70, 12).Value = 1.72
$ws.Cells.Item(70, 13).Value = "15/12/2023 17:50"
$ws.Cells.Item(70, 14).Value = 2.89
$ws.Cells.Item(70, 15).Value = "14/12/2023 01:12"
$ws.Cells.Item(70, 16).Value = 3.34
$ws.Cells.Item(70, 17).Value = "15/12/2023 17:55"
$ws.Cells.Item(70, 18).Value = 4.04
$ws.Cells.Item(70, 19).Value = "14/12/2023 01:12"
$ws.Cells.Item(70, 20).Value = 5.61
$ws.Cells.Item(70, 21).Value = "15/12/2023 17:17"
$ws.Cells.Item(70, 22).Value = "https://www.betexplorer.com/football/algeria/ligue-1/aso-chlef-paradou/2uPm8Fla/"

